# dados/teste_tamanho.xlsx — "Add files via upload" edit
#
# Net effect (derived from the OOXML diff):
#   1. Row 9's "contagem" (col D) value is corrected from 1 to 2.
#   2. The 18 rows with contagem = 0 and frase-count = 139/138 (old rows
#      20-37, the "neutral discourse" duplicates at the bottom of the
#      label=139/138 block) are deleted outright, shifting rows 38-46 up
#      to become rows 20-28.
#   3. An AutoFilter is (re)applied on column D ("contagem") of the
#      shrunken range A1:E28, showing only rows where contagem = 0 -
#      this hides the remaining rows 2-13 (contagem 1-3) without
#      deleting them.
#   4. The worksheet selection moves to B21, and the workbook's hidden
#      _FilterDatabase defined name is re-pointed at the new range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dados_completos")

# 1) Fix the mis-typed count for row 9 (id_deputado 178857).
$ws.Range("D9").Value2 = 2

# 2) Drop the 18 duplicate/neutral rows (old rows 20-37 inclusive).
#    Everything below shifts up to close the gap.
$ws.Range("A20:E37").EntireRow.Delete()

# Clear the pre-existing AutoFilter definition entirely first so that
# re-applying it below rebuilds it (and its <autoFilter ref=.../>) over
# the *new*, smaller A1:E28 range instead of reusing the stale A1:E46
# range left over from before the delete.
$ws.AutoFilterMode = $false

# 3) Re-apply AutoFilter over the new data extent, filtering column D
#    (4th column, 1-based) down to the discrete value "0" - this hides
#    rows 2-13 (contagem 1-3) while leaving rows 14-28 (contagem 0)
#    visible.
$criteria = @("0")
$ws.Range("A1:E28").AutoFilter(4, $criteria, 7)

# 4) Cosmetic bits that came along with the manual edit: selection
#    moved to B21, and the hidden filter-database defined name now
#    tracks the smaller range.
$ws.Range("B21").Select()

$name = $wb.Names.Item(1)
$name.RefersTo = "=Dados_completos!`$A`$1:`$E`$28"
